$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B10").Value = "field_wbddh_periodicity"
